$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix formatting of existing rows 7-8 (C:D) -----------------------------
# In the original file, C7:D8 use a "wrap text" style (cellXfs index 3).
# The edit folds that style into the existing "center / no-wrap" style
# (cellXfs index 2) that is already used by C2:D6 -- so we just re-apply the
# same horizontal alignment + wrap setting that those cells already carry.
$ws.Range("C7:D8").HorizontalAlignment = -4108
$ws.Range("C7:D8").WrapText = $false

# --- Append six new review rows (9-14) -------------------------------------
# Seed each new row by copying the formatting of an already-styled row
# (row 8) so the new cells land on the very same style indices instead of
# minting new ones, then overwrite with the real values.
$newRows = @(
  @("com.hamxa.shaynachim","bitcoin","mesikam455@gmail.com","imesika53@gmail.com","27/5/2019 15:59","incredible","yes"),
  @("com.hamxa.shaynachim","bitcoin","edenn0836@gmail.com","mesikam455@gmail.com","27/5/2019 15:59","The way I wanted to read it","yes"),
  @("com.hamxa.shaynachim","bitcoin","frimanoren6@gmail.com","edenn0836@gmail.com","27/5/2019 15:59","words of wisdom","yes"),
  @("com.hamxa.shaynachim","bitcoin","goldfinshmulik@gmail.com","frimanoren6@gmail.com","27/5/2019 15:59","clear and easy","yes"),
  @("com.hamxa.shaynachim","bitcoin","elad86643@gmail.com","goldfinshmulik@gmail.com","27/5/2019 15:59","great app","yes"),
  @("com.hamxa.shaynachim","bitcoin","erezadmoni26@gmail.com","goldfinshmulik@gmail.com","27/5/2019 15:59","wise and simple","yes")
)

$r = 9
foreach ($rowValues in $newRows) {
  $ws.Range("A8:G8").Copy($ws.Range("A" + $r + ":G" + $r))
  $ws.Rows.Item($r).RowHeight = 13.8

  $c = 1
  foreach ($val in $rowValues) {
    $ws.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }

  $r = $r + 1
}

# --- Keep the active selection in sync with the new last row ---------------
$null = $ws.Range("F15").Select()
